# Lab7BOM.xlsx update: "Added remaining footprints and aligned keys on layout"
#
# Adds four new BOM rows (152-155) to the "Parts" sheet describing the
# remaining components: a hot-swappable key socket (MISC), the LDO
# regulator and power-switch IC (POWER) and the TM4C123 ARM MCU (MCU).
# The running total in J4 (=SUM(J5:J200)) recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# --- Row 152: hot-swappable PCB socket (MISC / Amazon) ------------------
$ws.Cells.Item(152, 1).Value = 1
$ws.Cells.Item(152, 3).Value = "MISC"
$ws.Cells.Item(152, 4).Value = "Hot-swappable PCB Socket Sip Socket"
$ws.Cells.Item(152, 5).Value = "IJKT"
$ws.Cells.Item(152, 6).Value = "KS-2P02B01"
$ws.Cells.Item(152, 7).Value = "Amazon"
$ws.Cells.Item(152, 9).Value = 9.9
$ws.Cells.Item(152, 10).Formula = "=A152*I152"
$ws.Cells.Item(152, 12).Value = "https://www.amazon.com/dp/B0972HB9GY?ref=ppx_yo2ov_dt_b_product_details&th=1"

# --- Row 153: LDO voltage regulator (POWER / Mouser) ---------------------
$ws.Cells.Item(153, 1).Value = 1
$ws.Cells.Item(153, 3).Value = "POWER"
$ws.Cells.Item(153, 4).Value = "LDO Voltage Regulators Cap-Free NMOS 400mA"
$ws.Cells.Item(153, 5).Value = "Texas Instruments"
$ws.Cells.Item(153, 6).Value = "TPS73633DBV"
$ws.Cells.Item(153, 7).Value = "Mouser "
$ws.Cells.Item(153, 8).Value = "595-TPS73633DBVRG4"
$ws.Cells.Item(153, 9).Value = 1.5
$ws.Cells.Item(153, 10).Formula = "=A153*I153"
$ws.Cells.Item(153, 12).Value = "https://www.mouser.com/ProductDetail/Texas-Instruments/TPS73633DBVRG4?qs=6zVL%252ByCp0mpknSjwGe1Hbg%3D%3D"

# --- Row 154: power switch IC (POWER / Mouser) ----------------------------
$ws.Cells.Item(154, 1).Value = 1
$ws.Cells.Item(154, 3).Value = "POWER"
$ws.Cells.Item(154, 4).Value = "Power Switch IC"
$ws.Cells.Item(154, 5).Value = "Texas Instruments"
$ws.Cells.Item(154, 6).Value = "TPS2113ADRBR"
$ws.Cells.Item(154, 7).Value = "Mouser "
$ws.Cells.Item(154, 8).Value = "595-TPS2113ADRBR"
$ws.Cells.Item(154, 9).Value = 1.87
$ws.Cells.Item(154, 10).Formula = "=A154*I154"
$ws.Cells.Item(154, 12).Value = "https://www.mouser.com/ProductDetail/Texas-Instruments/TPS2113ADRBR?qs=g%2FrhRe7LVpRsXhRevikZ7Q%3D%3D"

# --- Row 155: TM4C123 ARM microcontroller (MCU / Mouser) ------------------
$ws.Cells.Item(155, 1).Value = 1
$ws.Cells.Item(155, 3).Value = "MCU"
$ws.Cells.Item(155, 4).Value = "ARM Microcontroller"
$ws.Cells.Item(155, 5).Value = "Texas Instruments"
$ws.Cells.Item(155, 6).Value = "TM4C123GH6PMI7"
$ws.Cells.Item(155, 7).Value = "Mouser "
$ws.Cells.Item(155, 8).Value = "595-TM4C123GH6PMI7"
$ws.Cells.Item(155, 9).Value = 12.5
$ws.Cells.Item(155, 10).Formula = "=A155*I155"
$ws.Cells.Item(155, 12).Value = "https://www.mouser.com/ProductDetail/Texas-Instruments/TM4C123GH6PMI7?qs=m%2F7bTylgptcEEuXHNU46tA%3D%3D"

# Keep the current selection on the sheet (author had landed on H76 before saving).
$ws.Range("H76").Select()
